# Apply "Add data for 2022-10-20" update to carjacking-by-month-yoy-latest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (the "Through" date label moves from 10-11 to 10-12)
$ws.Name = "Through 2022-10-12"

# Update the October row label to reflect new "through" date
$ws.Range("A11").Value = "October (through 10-12)"

# Update October row (row 11) values for columns B..I
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 19
$ws.Range("E11").Value = 30
$ws.Range("F11").Value = 14
$ws.Range("G11").Value = 55
$ws.Range("H11").Value = 80
$ws.Range("I11").Value = 39

# Update Total row (row 12) values for columns B..I
$ws.Range("B12").Value = 238
$ws.Range("C12").Value = 447
$ws.Range("D12").Value = 646
$ws.Range("E12").Value = 578
$ws.Range("F12").Value = 436
$ws.Range("G12").Value = 956
$ws.Range("H12").Value = 1327
$ws.Range("I12").Value = 1317
